# RF_07 analysis document - 20/08/2011
# "Realizar análise completa do requisito e materializar em documento"
#
# 1) The pré-condição about database access had a stray double space
#    ("...pleno acesso  ao banco...") that Word had flagged with a
#    gramStart/gramEnd proofing mark; normalize it to a single space.
# 2) Add a new pré-condição paragraph right after it, describing that the
#    app must contain registered muscles.

$d = $word.ActiveDocument

# --- Step 1: collapse the stray double space into a single space -----------
$d.Content.Find.Execute(
    "O aplicativo deve esta com pleno acesso  ao banco de dados", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "O aplicativo deve esta com pleno acesso ao banco de dados", 2)

# --- Step 2: locate that paragraph -----------------------------------------
$idx = 0
$targetIdx = -1
foreach ($para in $d.Paragraphs) {
    $idx = $idx + 1
    if ($para.Range.Text.Contains("O aplicativo deve esta com pleno acesso ao banco de dados")) {
        $targetIdx = $idx
    }
}

if ($targetIdx -eq -1) {
    throw "Could not locate the 'pleno acesso ao banco de dados' paragraph"
}

$target = $d.Paragraphs.Item($targetIdx)

# --- Step 3: append a new paragraph right after it, matching its style -----
$splitPoint = $target.Range.Duplicate
$splitPoint.Collapse(0)
$splitPoint.InsertAfter([char]13)

$newPara = $d.Paragraphs.Item($targetIdx + 1)
$newPara.Range.Text = "O aplicativo conter músculos cadastrados"
